# Add the 2022-Q3 quarterly data:
#  - insert a new "2022-Q3" worksheet (with per-fund holdings) right after
#    the "总计" (totals) sheet, before the existing "2022-Q2" sheet
#  - add a new summary row for 2022-Q3 at the top of the "总计" sheet,
#    pushing the existing 2022-Q2 / 2022-Q1 rows down

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)          # "总计"

# ---------------------------------------------------------------------
# 1. Shift the existing data rows on "总计" down by one, then write the
#    new 2022-Q3 summary row. We copy cell A3 -> A4 first so the new A4
#    cell inherits the bordered/centered index style, same for A2 -> A3.
# ---------------------------------------------------------------------
$b3 = $ws1.Range("B3").Value()
$c3 = $ws1.Range("C3").Value()
$d3 = $ws1.Range("D3").Value()

$b2 = $ws1.Range("B2").Value()
$c2 = $ws1.Range("C2").Value()
$d2 = $ws1.Range("D2").Value()

$ws1.Range("A3").Copy($ws1.Range("A4"))
$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = $b3
$ws1.Range("C4").Value = $c3
$ws1.Range("D4").Value = $d3

$ws1.Range("A2").Copy($ws1.Range("A3"))
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = $b2
$ws1.Range("C3").Value = $c2
$ws1.Range("D3").Value = $d2

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 0.54

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q3"

$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Clone the header row (B1:H1) and the index column (A2:A9) formatting
# from the existing "2022-Q2" sheet, BEFORE typing any values, so the new
# sheet matches the same bordered / bold / centered style already used on
# the other quarterly sheets.
$wsQ2.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$wsQ2.Range("A2").Copy($newSheet.Range("A2:A9"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("A3").Value = 1
$newSheet.Range("A4").Value = 2
$newSheet.Range("A5").Value = 3
$newSheet.Range("A6").Value = 4
$newSheet.Range("A7").Value = 5
$newSheet.Range("A8").Value = 6
$newSheet.Range("A9").Value = 7

# Columns B:G hold text values (fund codes / names / numeric-looking
# strings) in the source workbook, so force a text format before typing
# them in, then reset the style back to Normal so no stray styling is
# left behind on the cells.
$textRange = $newSheet.Range("B2:G9")
$textRange.NumberFormat = "@"

$newSheet.Range("B2").Value = "014269"
$newSheet.Range("C2").Value = "嘉实北交所精选两年定期混合A"
$newSheet.Range("D2").Value = "2.72"
$newSheet.Range("E2").Value = "90.37"
$newSheet.Range("F2").Value = "6.08"
$newSheet.Range("G2").Value = "0.1654"
$newSheet.Range("H2").Value = 7

$newSheet.Range("B3").Value = "014279"
$newSheet.Range("C3").Value = "汇添富北交所创新精选两年定开混合A"
$newSheet.Range("D3").Value = "3.20"
$newSheet.Range("E3").Value = "93.27"
$newSheet.Range("F3").Value = "4.69"
$newSheet.Range("G3").Value = "0.1501"
$newSheet.Range("H3").Value = 6

$newSheet.Range("B4").Value = "014663"
$newSheet.Range("C4").Value = "富国创新发展两年定期开放混合A"
$newSheet.Range("D4").Value = "2.24"
$newSheet.Range("E4").Value = "71.47"
$newSheet.Range("F4").Value = "3.50"
$newSheet.Range("G4").Value = "0.0784"
$newSheet.Range("H4").Value = 4

$newSheet.Range("B5").Value = "014271"
$newSheet.Range("C5").Value = "大成北交所两年定开混合A"
$newSheet.Range("D5").Value = "3.45"
$newSheet.Range("E5").Value = "65.31"
$newSheet.Range("F5").Value = "1.78"
$newSheet.Range("G5").Value = "0.0614"
$newSheet.Range("H5").Value = 9

$newSheet.Range("B6").Value = "014270"
$newSheet.Range("C6").Value = "嘉实北交所精选两年定期混合C"
$newSheet.Range("D6").Value = "0.53"
$newSheet.Range("E6").Value = "90.37"
$newSheet.Range("F6").Value = "6.08"
$newSheet.Range("G6").Value = "0.0322"
$newSheet.Range("H6").Value = 7

$newSheet.Range("B7").Value = "014280"
$newSheet.Range("C7").Value = "汇添富北交所创新精选两年定开混合C"
$newSheet.Range("D7").Value = "0.51"
$newSheet.Range("E7").Value = "93.27"
$newSheet.Range("F7").Value = "4.69"
$newSheet.Range("G7").Value = "0.0239"
$newSheet.Range("H7").Value = 6

$newSheet.Range("B8").Value = "014272"
$newSheet.Range("C8").Value = "大成北交所两年定开混合C"
$newSheet.Range("D8").Value = "0.82"
$newSheet.Range("E8").Value = "65.31"
$newSheet.Range("F8").Value = "1.78"
$newSheet.Range("G8").Value = "0.0146"
$newSheet.Range("H8").Value = 9

$newSheet.Range("B9").Value = "014664"
$newSheet.Range("C9").Value = "富国创新发展两年定期开放混合C"
$newSheet.Range("D9").Value = "0.33"
$newSheet.Range("E9").Value = "71.47"
$newSheet.Range("F9").Value = "3.50"
$newSheet.Range("G9").Value = "0.0116"
$newSheet.Range("H9").Value = 4

$textRange.Style = "Normal"

# ---------------------------------------------------------------------
# 3. Restore "总计" as the active sheet (it was the active sheet before
#    this edit, and adding a worksheet normally activates the new one).
# ---------------------------------------------------------------------
$ws1.Activate()
